$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.982.28"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.418.35"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.83"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.20"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  +1.42%  "
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.21"
$ws.Range("E11").Value = "  -2.98%  "
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.77"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "2.853.95"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "61.955.23"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "2.422.07"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.25"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "322.77"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.85"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.14"
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.55"
$ws.Range("E23").Value = "  +2.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.75"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.80"
$ws.Range("E25").Value = "  -1.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "562.93"
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.537.57"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").Value = "0.0₃0941"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("E31").Value = "  -3.86%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.51"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.77"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "155.30"
$ws.Range("E38").Value = "  +3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.44"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.54"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.82"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.28"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.92"
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.86"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.593"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0922"
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0227"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +0.69%  "
